$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout:
# Row1: B1=0
# Row2: A2="Years at Company", B2=1
# Row3: A3="Monthly Income", B3=1226
# Row4: A4="Number of Promotions", B4=0
# Row5: A5="Distance from Home", B5=1
# Row6: A6="Number of Dependents", B6=0
#
# Target layout:
# Row1: B1=0
# Row2: A2="Years at Company", B2=1
# Row3: A3="Monthly Income", B3=1226
# Row4: A4="Distance from Home", B4=1

# Delete row 4 ("Number of Promotions") entirely, shifting rows up.
$ws.Rows(4).Delete()

# Delete the (now) last row which held "Number of Dependents".
$ws.Rows(5).Delete()
